$wb = $excel.ActiveWorkbook

# --- Sheet "YDS": extend Week-13 yardage logs ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '8 2 6 9 1 8 -3 1 3 1 3 3 3 7 3 2 2 11 2 6 15 25 2 3 46 4 6 1 10 1 -1 2 5 4 0 8 4 4 4 1 4 15 2 5 0 0 7 3 4 2 8 1 0 4 6 8 1 -2 4 4 8 2 3 4 3 2 2 0 3 0 7 6 5 1 6 9 1 2 2 3 4 3 1 8 0 17 8 -5 7 9 6 3 0 8 5 41 12 7 7 1 5 1 4 3 8 2 5 -1 4 0 8 22 2 11 9 9 8 2 -1 3 5 -1 0 0 3 1 4 12 4 6 6 0 8 3 1 3 14 3 7 3 -3 7 5 2 1 7 5 2 3 4 0 4 5 0 4 5 4 -1 6 34 1 0 -1 -2 6 5 7 13 2 6 2 12 3 -1 4 0 7 15 12 -3 1 22 4 -1 1 7 5 4 2 3 -1 15 0 -4 4 0 4 11 7 13 24 8 15 0 1 7 6 3 9 2 3 10 0 3 8 9 4 28 3 8 5 2 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1'
$wsYDS.Range("C2").Value = '2 4 2 2 -4 2 0 9 2 1 2 25 -2 5 18 3 8 -1 6 18 -1 2 -1 6 2 1 3 2 8 3 4 -2 6 1 12 1 8 2 4 2 7 -1 4 1 -1 11 4 -4 7 3 1 7 4 7 6 0 1 0 5 1 10 3 4 0 12 0 -6 2 -1 2 4 1 3 2 4 2 5 4 2 10 5 2 3 15 3 6 4 1 9 5 0 -1 23 9 -1 -2 9 5 9 0 -1 2 4 3 76 2 -3 7 4 6 3 19 0 4 3 1 2 3 0 -1 13 1 7 5 5 5 -3 2 2 5 3 -2 12 3 1 1 0 2 6 2 0 3 1 7 6 3 2 3 -4 -1 15 1 -2 3 0 4 11 7 -1 5 1 6 2 8 2 0 -2 3 5 5 15 6 5 4 -4 5 -2 -2 5 8 5 13 0 0 0 5 1 -1 3 8 10 6 10 3 3 11 2 4 5 1 6 0 18 4 3 7 -1 2 2 3 5 40 10 9 0 3 2 1 19 6 5 18 9 3 -3 1 5 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2'
$wsYDS.Range("B3").Value = '7 -2 6 5 10 16 13 4 7 12 10 5 37 10 8 3 7 9 9 12 7 4 6 12 -1 13 10 13 5 1 35 5 6 1 7 7 22 13 6 7 8 6 9 2 41 3 23 11 13 1 28 3 7 9 7 17 4 7 14 18 2 41 7 13 22 11 5 11 9 5 6 13 10 11 6 5 19 7 5 25 31 9 20 1 16 7 8 5 16 21 2 37 9 5 5 1 19 23 8 24 35 61 15 53 41 1 -2 11 5 16 16 8 9 12 12 6 4 3 11 6 7 3 3 21 19 2 14 5 17 12 29 1 5 9 11 18 9 15 3 4 1 6 6 31 12 7 20 1 11 8 1 20 3 15 7 8 6 15 4 10 6 3 3 8 8 17 11 14 19 13 1 11 -1 20 8 4 3 7 1 16 0 4 9 9 14 1 7 2 9 11 13 7 9 5 22 8 12 15 3 5 -2 28 1 7 19 5 14 3 11 17 15 4 7 7 12 3 6 8 6 6 16 57 12 23 31 49 25 8 36 20 43 7 11 9 7 20 2 18 3 6 14 12 7 6 10 6 7 16 7 31 1 9 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6'
$wsYDS.Range("C3").Value = '14 6 6 12 9 10 6 19 22 5 9 24 4 -2 5 24 14 1 13 7 6 9 5 18 12 4 5 14 4 11 16 11 3 5 8 0 4 2 7 8 4 5 6 0 11 17 10 73 14 7 4 15 1 11 15 37 3 8 7 8 5 8 12 3 6 12 12 6 15 11 9 7 8 5 8 10 6 26 8 12 8 11 5 4 3 11 8 17 6 8 5 5 4 1 8 6 4 4 11 11 5 48 7 8 8 8 4 13 20 6 6 14 22 5 8 9 23 7 8 0 14 15 22 12 8 3 8 5 1 -1 6 12 10 4 10 -2 11 4 6 12 -1 6 7 5 5 4 5 28 9 4 16 9 4 6 7 3 3 2 9 10 4 21 18 14 28 19 14 12 11 11 6 9 -7 21 18 8 -4 18 18 5 23 9 18 0 3 8 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 12 7'

# --- Sheet "OFF": weekly offensive totals ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 110
$wsOFF.Range("E2").Value = 7
$wsOFF.Range("F2").Value = 54
$wsOFF.Range("G2").Value = 39
$wsOFF.Range("I2").Value = 7
$wsOFF.Range("J2").Value = 24
$wsOFF.Range("L2").Value = 233
$wsOFF.Range("M2").Value = 148
$wsOFF.Range("O2").Value = 12
$wsOFF.Range("Q2").Value = 384
$wsOFF.Range("C3").Value = 170
$wsOFF.Range("E3").Value = 22
$wsOFF.Range("F3").Value = 77
$wsOFF.Range("G3").Value = 22
$wsOFF.Range("H3").Value = 24
$wsOFF.Range("I3").Value = 30
$wsOFF.Range("J3").Value = 47
$wsOFF.Range("N3").Value = 7

# --- Sheet "DEF": weekly defensive totals ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 133
$wsDEF.Range("E2").Value = 11
$wsDEF.Range("F2").Value = 43
$wsDEF.Range("G2").Value = 39
$wsDEF.Range("H2").Value = 5
$wsDEF.Range("I2").Value = 7
$wsDEF.Range("J2").Value = 22
$wsDEF.Range("L2").Value = 140
$wsDEF.Range("M2").Value = 77
$wsDEF.Range("O2").Value = 19
$wsDEF.Range("P2").Value = 9
$wsDEF.Range("Q2").Value = 316
$wsDEF.Range("C3").Value = 106
$wsDEF.Range("E3").Value = 36
$wsDEF.Range("H3").Value = 30

# --- Sheet "ST": special teams totals + extended logs ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 69
$wsST.Range("D2").Value = 33
$wsST.Range("F2").Value = 94
$wsST.Range("G2").Value = 92
$wsST.Range("J2").Value = 36
$wsST.Range("K2").Value = 32
$wsST.Range("B3").Value = 38
$wsST.Range("D3").Value = '55 46 36 42 41 54 37 60 44 47 47 48 45 32 35 53 40 55 19 41 48 59 41 48 41 37 55 48 43 49 33 34 59'
$wsST.Range("D4").Value = '0 5 0 4 0 0 0 0 6 30 30 0 0 0 0 2 0 8 0 -4 8 14 0 0 8 0 7 12 0 0 0 0 0'
$wsST.Range("D5").Value = '0 1 0 0 0 20 0 0 0 13 6 33 0 0 5 0 9 6 12 -1 7 0 6 14 16 0 9 0 0 0 0 0 0 0 26 0 0 11 10 1 18 8 0 0 0 0 0'
$wsST.Range("B6").Value = '75 23 19 26 0 20 0 30 18 19 18 45 33 28 23 22 12 12 27 20 17 17 27'

# --- Sheet "TURNS": weekly turnover totals ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("D2").Value = 9
$wsTURNS.Range("E2").Value = 9

# --- Sheet "PEN": weekly penalty totals ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 15
$wsPEN.Range("B3").Value = 19

